$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previous template contents entirely (old layout was A1:O4).
$ws.Range("A1:O10").ClearContents()

# New header row.
$ws.Range("A1").Value = "委托时间"
$ws.Range("B1").Value = "证券代码"
$ws.Range("C1").Value = "证券名称"
$ws.Range("D1").Value = "操作"
$ws.Range("E1").Value = "备注"
$ws.Range("F1").Value = "委托数量"
$ws.Range("G1").Value = "成交数量"
$ws.Range("H1").Value = "撤消数量"
$ws.Range("I1").Value = "成交金额"
$ws.Range("J1").Value = "委托价格"
$ws.Range("K1").Value = "成交均价"
$ws.Range("L1").Value = "合同编号"
$ws.Range("M1").Value = "交易类别"

# New data rows (2-8).
$ws.Range("A2").Value = 0.43079861111111112
$ws.Range("B2").Value = 2798
$ws.Range("C2").Value = "帝王洁具"
$ws.Range("D2").Value = "证券买入"
$ws.Range("E2").Value = "已成"
$ws.Range("F2").Value = 4100
$ws.Range("G2").Value = 4100
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 152900
$ws.Range("J2").Value = 37.299999999999997
$ws.Range("K2").Value = 37.292999999999999
$ws.Range("L2").Value = 56987
$ws.Range("M2").Value = "波段"
$ws.Range("A3").Value = 0.43062500000000004
$ws.Range("B3").Value = 2798
$ws.Range("C3").Value = "帝王洁具"
$ws.Range("D3").Value = "证券买入"
$ws.Range("E3").Value = "部撤"
$ws.Range("F3").Value = 8100
$ws.Range("G3").Value = 2600
$ws.Range("H3").Value = 5500
$ws.Range("I3").Value = 96555
$ws.Range("J3").Value = 37.15
$ws.Range("K3").Value = 37.137
$ws.Range("L3").Value = 56809
$ws.Range("M3").Value = "波段"
$ws.Range("A4").Value = 0.43056712962962962
$ws.Range("B4").Value = 2798
$ws.Range("C4").Value = "帝王洁具"
$ws.Range("D4").Value = "证券买入"
$ws.Range("E4").Value = "已成"
$ws.Range("F4").Value = 5100
$ws.Range("G4").Value = 5100
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 189312
$ws.Range("J4").Value = 37.119999999999997
$ws.Range("K4").Value = 37.119999999999997
$ws.Range("L4").Value = 56755
$ws.Range("M4").Value = "波段"
$ws.Range("A5").Value = 0.4305208333333333
$ws.Range("B5").Value = 2798
$ws.Range("C5").Value = "帝王洁具"
$ws.Range("D5").Value = "证券买入"
$ws.Range("E5").Value = "部撤"
$ws.Range("F5").Value = 4100
$ws.Range("G5").Value = 3500
$ws.Range("H5").Value = 600
$ws.Range("I5").Value = 129885
$ws.Range("J5").Value = 37.11
$ws.Range("K5").Value = 37.11
$ws.Range("L5").Value = 56712
$ws.Range("M5").Value = "波段"
$ws.Range("A6").Value = 0.43048611111111112
$ws.Range("B6").Value = 2798
$ws.Range("C6").Value = "帝王洁具"
$ws.Range("D6").Value = "证券买入"
$ws.Range("E6").Value = "已成"
$ws.Range("F6").Value = 2100
$ws.Range("G6").Value = 2100
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 77931
$ws.Range("J6").Value = 37.11
$ws.Range("K6").Value = 37.11
$ws.Range("L6").Value = 56676
$ws.Range("M6").Value = "波段"
$ws.Range("A7").Value = 0.4304398148148148
$ws.Range("B7").Value = 2798
$ws.Range("C7").Value = "帝王洁具"
$ws.Range("D7").Value = "证券买入"
$ws.Range("E7").Value = "已成"
$ws.Range("F7").Value = 1400
$ws.Range("G7").Value = 1400
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 51954
$ws.Range("J7").Value = 37.11
$ws.Range("K7").Value = 37.11
$ws.Range("L7").Value = 56637
$ws.Range("M7").Value = "波段"
$ws.Range("A8").Value = 0.43089120370370365
$ws.Range("B8").Value = 2798
$ws.Range("C8").Value = "帝王洁具"
$ws.Range("D8").Value = "证券买入"
$ws.Range("E8").Value = "已成"
$ws.Range("F8").Value = 5100
$ws.Range("G8").Value = 5100
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 190725
$ws.Range("J8").Value = 37.479999999999997
$ws.Range("K8").Value = 37.396999999999998
$ws.Range("L8").Value = 57070
$ws.Range("M8").Value = "波段"

# Column A (委托时间) holds time-of-day serials; format as h:mm:ss.
$ws.Range("A2:A8").NumberFormat = "h:mm:ss"

# Restore selection as captured after the edit.
$ws.Range("Q21").Select()
